$d = $word.ActiveDocument
$d.Content.Delete()
